# MedicationReport.xlsx — replace "brand"/"generic" record placeholders with
# "title"/"subtitle" (AB#15894), tidy the two-row sample layout, and drop the
# now-unused border formatting from the sample header/row cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sample data cells: swap the Brand/Generic columns for Title/Subtitle ---
# Column C used to carry "{d.records[...].brand}", column D carried
# "{d.records[...].generic}" — the model now exposes title/subtitle instead.
$ws.Range("C2").Value = "{d.records[i].title}"
$ws.Range("C3").Value = "{d.records[i+1].title}"
$ws.Range("D2").Value = "{d.records[i].subtitle}"
$ws.Range("D3").Value = "{d.records[i+1].subtitle}"

# --- Drop the (empty/invisible) border formatting that's no longer applied ---
$ws.Range("A1").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("A2").Borders.LineStyle = -4142
$ws.Range("A3").Borders.LineStyle = -4142

# --- Column widths (slightly narrower, matching the refreshed template) ---
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 34.833333333333336
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666
$ws.Columns.Item(8).ColumnWidth = 29.833333333333332
$ws.Columns.Item(9).ColumnWidth = 43.5

# --- Row height for the first sample row ---
$ws.Rows.Item(2).RowHeight = 16.15

# --- Reset the lingering selection left over from editing, back to A1 ---
$ws.Range("A1").Select()

# --- Page orientation (now explicit in the saved template) ---
$ws.PageSetup.Orientation = 1   # xlPortrait
